$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "顺灏股份"
$ws.Cells.Item(2, 2).Value = "平潭发展"
$ws.Cells.Item(2, 3).Value = "平潭发展"
$ws.Cells.Item(3, 1).Value = "美年健康"
$ws.Cells.Item(3, 2).Value = "中金公司"
$ws.Cells.Item(3, 3).Value = "合富中国"
$ws.Cells.Item(4, 1).Value = "航天发展"
$ws.Cells.Item(4, 2).Value = "永辉超市"
$ws.Cells.Item(4, 3).Value = "百大集团"
$ws.Cells.Item(5, 1).Value = "平潭发展"
$ws.Cells.Item(5, 2).Value = "C沐曦-U"
$ws.Cells.Item(5, 3).Value = "航天发展"
$ws.Cells.Item(6, 1).Value = "西部材料"
$ws.Cells.Item(6, 2).Value = "顺灏股份"
$ws.Cells.Item(6, 3).Value = "永辉超市"
$ws.Cells.Item(7, 1).Value = "永辉超市"
$ws.Cells.Item(7, 2).Value = "航天发展"
$ws.Cells.Item(7, 3).Value = "顺灏股份"
$ws.Cells.Item(8, 1).Value = "航天电子"
$ws.Cells.Item(8, 2).Value = "中国卫星"
$ws.Cells.Item(8, 3).Value = "浙江世宝"
$ws.Cells.Item(9, 1).Value = "中国卫星"
$ws.Cells.Item(9, 2).Value = "航天电子"
$ws.Cells.Item(9, 3).Value = "美年健康"
$ws.Cells.Item(10, 1).Value = "百大集团"
$ws.Cells.Item(10, 2).Value = "东兴证券"
$ws.Cells.Item(10, 3).Value = "再升科技"
$ws.Cells.Item(11, 1).Value = "合富中国"
$ws.Cells.Item(11, 2).Value = "合富中国"
$ws.Cells.Item(11, 3).Value = "东百集团"
$ws.Cells.Item(12, 1).Value = "天银机电"
$ws.Cells.Item(12, 2).Value = "百大集团"
$ws.Cells.Item(12, 3).Value = "航天电子"
$ws.Cells.Item(13, 1).Value = "浙江世宝"
$ws.Cells.Item(13, 2).Value = "西部材料"
$ws.Cells.Item(13, 3).Value = "西部材料"
$ws.Cells.Item(14, 1).Value = "中金公司"
$ws.Cells.Item(14, 2).Value = "天银机电"
$ws.Cells.Item(14, 3).Value = "雪人集团"
$ws.Cells.Item(15, 1).Value = "鹭燕医药"
$ws.Cells.Item(15, 2).Value = "美年健康"
$ws.Cells.Item(15, 3).Value = "华菱线缆"
$ws.Cells.Item(16, 1).Value = "C沐曦-U"
$ws.Cells.Item(16, 2).Value = "龙洲股份"
$ws.Cells.Item(16, 3).Value = "利群股份"
$ws.Cells.Item(17, 1).Value = "英维克"
$ws.Cells.Item(17, 2).Value = "浙江世宝"
$ws.Cells.Item(17, 3).Value = "博纳影业"
$ws.Cells.Item(18, 1).Value = "东百集团"
$ws.Cells.Item(18, 2).Value = "飞龙股份"
$ws.Cells.Item(18, 3).Value = "中央商场"
$ws.Cells.Item(19, 1).Value = "龙洲股份"
$ws.Cells.Item(19, 2).Value = "英维克"
$ws.Cells.Item(19, 3).Value = "航天机电"
$ws.Cells.Item(20, 1).Value = "南京商旅"
$ws.Cells.Item(20, 2).Value = "海峡创新"
$ws.Cells.Item(20, 3).Value = "飞龙股份"
$ws.Cells.Item(21, 1).Value = "飞龙股份"
$ws.Cells.Item(21, 2).Value = "利群股份"
$ws.Cells.Item(21, 3).Value = "龙洲股份"
